# upgrade left table until javakheti
# Zestaponi.xlsx: fix municipality name typo (Zestafoni -> Zestaponi) and
# extend the employees table with the 2023 column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Zestafoni" -> "Zestaponi" typo -------------------------------
$ws.Name = "Zestaponi"
$ws.Range("A1").Value = "Number of employees of business sector in Zestaponi Municipality"

# --- Clean up the stray bottom border under the year header row ------------
# (B3:J3 previously closed off with a bottom border; the table now continues
# past column J, so only the top border of that row should remain.)
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142

# --- Extend the table with the new 2023 column (K) --------------------------
# Clone the formatting of the last "clean" data column (I) onto the new
# column K, then fill in the 2023 figures.
$ws.Range("I3:I6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 6455
$ws.Range("K5").Value = 1589
$ws.Range("K6").Value = 4866

# Close off the right edge of the table with a thin border on the new column.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2
